# ---------------------------------------------------------------------------
# Scheduled market-data refresh for the Belias server Leve-profit workbook.
#
# Each worksheet (one per crafting class: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# lists Leves together with live marketboard pricing pulled in columns H:N:
#   H currentAveragePrice        I currentAveragePriceNQ  J currentAveragePriceHQ
#   K LevePriceNQ                L LevePriceHQ
#   M LeveProfitNQ               N LeveProfitHQ
#
# This runner simply re-stamps the refreshed price/profit figures for the
# rows whose market data changed since the last run. A few rows also lose
# their HQ-profit figure (column N) entirely when no HQ listing exists any
# more, matching upstream (ClearContents, not just blanking the value).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: Distill, My Heart / Distilled Water
$ws.Cells.Item(9, 8).Value = 1985
$ws.Cells.Item(9, 9).Value = 600
$ws.Cells.Item(9, 10).Value = 2578.5715
$ws.Cells.Item(9, 11).Value = 600
$ws.Cells.Item(9, 12).Value = 2578.5715
$ws.Cells.Item(9, 13).Value = -431
$ws.Cells.Item(9, 14).Value = -2916.5715

# Row 42: Eye of the Beholder / Hi-Potion of Dexterity
$ws.Cells.Item(42, 8).Value = 86.25
$ws.Cells.Item(42, 9).Value = 48.333332
$ws.Cells.Item(42, 10).Value = 200
$ws.Cells.Item(42, 11).Value = 144.999996
$ws.Cells.Item(42, 12).Value = 600
$ws.Cells.Item(42, 13).Value = 85.00000399999999
$ws.Cells.Item(42, 14).Value = -1060

# Row 116: Growing Up / Growth Formula Kappa
$ws.Cells.Item(116, 8).Value = 2808.75
$ws.Cells.Item(116, 9).Value = 2745
$ws.Cells.Item(116, 10).Value = 3000
$ws.Cells.Item(116, 11).Value = 2745
$ws.Cells.Item(116, 12).Value = 3000
$ws.Cells.Item(116, 13).Value = 697
$ws.Cells.Item(116, 14).Value = -9884

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 1497.8445
$ws.Cells.Item(137, 9).Value = 1371.24
$ws.Cells.Item(137, 10).Value = 1656.1
$ws.Cells.Item(137, 11).Value = 4113.72
$ws.Cells.Item(137, 12).Value = 4968.299999999999
$ws.Cells.Item(137, 13).Value = -1563.72
$ws.Cells.Item(137, 14).Value = -10068.3

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 1949.1613
$ws.Cells.Item(138, 9).Value = 1264.5686
$ws.Cells.Item(138, 10).Value = 2780.4524
$ws.Cells.Item(138, 11).Value = 3793.7058
$ws.Cells.Item(138, 12).Value = 8341.3572
$ws.Cells.Item(138, 13).Value = 1346.2942
$ws.Cells.Item(138, 14).Value = -18621.3572

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 2195.5789
$ws.Cells.Item(141, 9).Value = 1190.4642
$ws.Cells.Item(141, 10).Value = 5009.9
$ws.Cells.Item(141, 11).Value = 3571.3926
$ws.Cells.Item(141, 12).Value = 15029.7
$ws.Cells.Item(141, 13).Value = 1608.6074
$ws.Cells.Item(141, 14).Value = -25389.7

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value = 4600
$ws.Cells.Item(2, 9).Value = 4600
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 4600
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -4487
$ws.Cells.Item(2, 14).ClearContents()

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Cells.Item(88, 8).Value = 1699
$ws.Cells.Item(88, 9).Value = 1699
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 1699
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = -1293
$ws.Cells.Item(88, 14).ClearContents()

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Cells.Item(91, 8).Value = 1699
$ws.Cells.Item(91, 9).Value = 1699
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 1699
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = -295
$ws.Cells.Item(91, 14).ClearContents()

# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 4600
$ws.Cells.Item(116, 9).Value = 4600
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 4600
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -2306
$ws.Cells.Item(116, 14).ClearContents()

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 1476.6957
$ws.Cells.Item(132, 9).Value = 1172.6666
$ws.Cells.Item(132, 10).Value = 2571.2
$ws.Cells.Item(132, 11).Value = 3517.9998
$ws.Cells.Item(132, 12).Value = 7713.599999999999
$ws.Cells.Item(132, 13).Value = -987.9998000000001
$ws.Cells.Item(132, 14).Value = -12773.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value = 4600
$ws.Cells.Item(3, 9).Value = 4600
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 4600
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -4486
$ws.Cells.Item(3, 14).ClearContents()

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 934901.25
$ws.Cells.Item(134, 9).Value = 1433264.5
$ws.Cells.Item(134, 10).Value = 4623.2
$ws.Cells.Item(134, 11).Value = 4299793.5
$ws.Cells.Item(134, 12).Value = 13869.6
$ws.Cells.Item(134, 13).Value = -4297258.5
$ws.Cells.Item(134, 14).Value = -18939.6

$ws = $wb.Worksheets.Item("CRP")
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Cells.Item(62, 8).Value = 2282.353
$ws.Cells.Item(62, 9).Value = 2383.3333
$ws.Cells.Item(62, 10).Value = 2040
$ws.Cells.Item(62, 11).Value = 2383.3333
$ws.Cells.Item(62, 12).Value = 2040
$ws.Cells.Item(62, 13).Value = -1759.3333
$ws.Cells.Item(62, 14).Value = -3288

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Cells.Item(65, 8).Value = 2282.353
$ws.Cells.Item(65, 9).Value = 2383.3333
$ws.Cells.Item(65, 10).Value = 2040
$ws.Cells.Item(65, 11).Value = 11916.6665
$ws.Cells.Item(65, 12).Value = 10200
$ws.Cells.Item(65, 13).Value = -8796.666499999999
$ws.Cells.Item(65, 14).Value = -16440

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Cells.Item(122, 8).Value = 1184.8334
$ws.Cells.Item(122, 10).Value = 1215.6666
$ws.Cells.Item(122, 12).Value = 3646.9998
$ws.Cells.Item(122, 14).Value = -8546.9998

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Cells.Item(134, 8).Value = 5155.636
$ws.Cells.Item(134, 9).Value = 4489
$ws.Cells.Item(134, 10).Value = 6933.3335
$ws.Cells.Item(134, 11).Value = 13467
$ws.Cells.Item(134, 12).Value = 20800.0005
$ws.Cells.Item(134, 13).Value = -10932
$ws.Cells.Item(134, 14).Value = -25870.0005

# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Cells.Item(141, 8).Value = 92919.91
$ws.Cells.Item(141, 10).Value = 92919.91
$ws.Cells.Item(141, 12).Value = 92919.91
$ws.Cells.Item(141, 14).Value = -103279.91

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Cells.Item(5, 8).Value = 663.8605
$ws.Cells.Item(5, 9).Value = 306.23077
$ws.Cells.Item(5, 11).Value = 918.69231
$ws.Cells.Item(5, 13).Value = -806.69231

# Row 101: No Othard Choice / Egg Foo Young
$ws.Cells.Item(101, 8).Value = 8000
$ws.Cells.Item(101, 10).Value = 8000
$ws.Cells.Item(101, 12).Value = 24000
$ws.Cells.Item(101, 14).Value = -28868

# Row 107: Slippery Service / Frantoio Oil
$ws.Cells.Item(107, 8).Value = 764
$ws.Cells.Item(107, 9).Value = 844.2857
$ws.Cells.Item(107, 10).Value = 202
$ws.Cells.Item(107, 11).Value = 2532.8571
$ws.Cells.Item(107, 12).Value = 606
$ws.Cells.Item(107, 13).Value = -612.8571000000002
$ws.Cells.Item(107, 14).Value = -4446

# Row 129: Comfort Food / Yakow Moussaka
$ws.Cells.Item(129, 8).Value = 1859.1082
$ws.Cells.Item(129, 9).Value = 1048.8235
$ws.Cells.Item(129, 10).Value = 2547.85
$ws.Cells.Item(129, 11).Value = 3146.4705
$ws.Cells.Item(129, 12).Value = 7643.549999999999
$ws.Cells.Item(129, 13).Value = 1853.5295
$ws.Cells.Item(129, 14).Value = -17643.55

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 951.05
$ws.Cells.Item(131, 10).Value = 983.6421
$ws.Cells.Item(131, 12).Value = 2950.9263
$ws.Cells.Item(131, 14).Value = -13030.9263

# Row 132: More Mezcal / Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 2323.3333
$ws.Cells.Item(132, 9).Value = 2644.2856
$ws.Cells.Item(132, 10).Value = 1200
$ws.Cells.Item(132, 11).Value = 23798.5704
$ws.Cells.Item(132, 12).Value = 10800
$ws.Cells.Item(132, 13).Value = -21268.5704
$ws.Cells.Item(132, 14).Value = -15860

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Cells.Item(135, 8).Value = 663.8605
$ws.Cells.Item(135, 9).Value = 306.23077
$ws.Cells.Item(135, 11).Value = 2756.07693
$ws.Cells.Item(135, 13).Value = -221.0769300000002

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Cells.Item(80, 8).Value = 2327.3845
$ws.Cells.Item(80, 9).Value = 2349.875
$ws.Cells.Item(80, 10).Value = 2291.4
$ws.Cells.Item(80, 11).Value = 2349.875
$ws.Cells.Item(80, 12).Value = 2291.4
$ws.Cells.Item(80, 13).Value = -1351.875
$ws.Cells.Item(80, 14).Value = -4287.4

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Cells.Item(83, 8).Value = 2327.3845
$ws.Cells.Item(83, 9).Value = 2349.875
$ws.Cells.Item(83, 10).Value = 2291.4
$ws.Cells.Item(83, 11).Value = 11749.375
$ws.Cells.Item(83, 12).Value = 11457
$ws.Cells.Item(83, 13).Value = -6757.375
$ws.Cells.Item(83, 14).Value = -21441

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Cells.Item(102, 8).Value = 1203.5
$ws.Cells.Item(102, 9).Value = 1180.1538
$ws.Cells.Item(102, 11).Value = 1180.1538
$ws.Cells.Item(102, 13).Value = 441.8462

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Cells.Item(122, 8).Value = 988.44446
$ws.Cells.Item(122, 9).Value = 918
$ws.Cells.Item(122, 10).Value = 1076.5
$ws.Cells.Item(122, 11).Value = 2754
$ws.Cells.Item(122, 12).Value = 3229.5
$ws.Cells.Item(122, 13).Value = -304
$ws.Cells.Item(122, 14).Value = -8129.5

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 2114.6155
$ws.Cells.Item(126, 9).Value = 1449
$ws.Cells.Item(126, 10).Value = 4333.3335
$ws.Cells.Item(126, 11).Value = 4347
$ws.Cells.Item(126, 12).Value = 13000.0005
$ws.Cells.Item(126, 13).Value = -1877
$ws.Cells.Item(126, 14).Value = -17940.0005

# Row 132: On Board for Lar / Lar Ingot
$ws.Cells.Item(132, 8).Value = 3656.0356
$ws.Cells.Item(132, 9).Value = 2235.862
$ws.Cells.Item(132, 10).Value = 5181.407
$ws.Cells.Item(132, 11).Value = 6707.586
$ws.Cells.Item(132, 12).Value = 15544.221
$ws.Cells.Item(132, 13).Value = -4177.586
$ws.Cells.Item(132, 14).Value = -20604.221

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value = 3384.3572
$ws.Cells.Item(40, 9).Value = 3686.889
$ws.Cells.Item(40, 10).Value = 2839.8
$ws.Cells.Item(40, 11).Value = 3686.889
$ws.Cells.Item(40, 12).Value = 2839.8
$ws.Cells.Item(40, 13).Value = -3550.889
$ws.Cells.Item(40, 14).Value = -3111.8

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Cells.Item(68, 8).Value = 15722
$ws.Cells.Item(68, 10).Value = 1725
$ws.Cells.Item(68, 12).Value = 1725
$ws.Cells.Item(68, 14).Value = -3223

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Cells.Item(71, 8).Value = 15722
$ws.Cells.Item(71, 10).Value = 1725
$ws.Cells.Item(71, 12).Value = 8625
$ws.Cells.Item(71, 14).Value = -16113

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Cells.Item(113, 8).Value = 617
$ws.Cells.Item(113, 9).Value = 385.66666
$ws.Cells.Item(113, 10).Value = 1137.5
$ws.Cells.Item(113, 11).Value = 1156.99998
$ws.Cells.Item(113, 12).Value = 3412.5
$ws.Cells.Item(113, 13).Value = 1013.00002
$ws.Cells.Item(113, 14).Value = -7752.5

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 950853.5
$ws.Cells.Item(132, 9).Value = 1402450.4
$ws.Cells.Item(132, 10).Value = 2500.15
$ws.Cells.Item(132, 11).Value = 4207351.199999999
$ws.Cells.Item(132, 12).Value = 7500.450000000001
$ws.Cells.Item(132, 13).Value = -4204821.199999999
$ws.Cells.Item(132, 14).Value = -12560.45

